# Insert a new data row at row 54 (pushing the existing rows 54-102 down to
# 55-103), then fill the newly inserted row with its own values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(54).Insert()

$ws.Range("A54").Value2 = 4
$ws.Range("B54").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C54").Value2 = "Los Lagos"
$ws.Range("D54").Value2 = 44827
$ws.Range("E54").Value2 = 10
$ws.Range("F54").Value2 = 100112026
$ws.Range("G54").Value2 = "Haba"
$ws.Range("H54").Value2 = "Sin especificar"
$ws.Range("I54").Value2 = "Primera"
$ws.Range("J54").Value2 = 80
$ws.Range("K54").Value2 = 13000
$ws.Range("L54").Value2 = 13000
$ws.Range("M54").Value2 = 13000
$ws.Range("N54").Value2 = "`$/saco 25 kilos"
$ws.Range("O54").Value2 = "Región Metropolitana"
$ws.Range("P54").Value2 = 520
$ws.Range("Q54").Value2 = 25
$ws.Range("R54").Value2 = "Hortaliza"
